$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''245.95'
$ws.Range("G2").Value = '''13'
$ws.Range("D3").Value = '''26.17'
$ws.Range("G3").Value = '''13'
$ws.Range("D4").Value = '''5.075'
$ws.Range("G4").Value = '''13'
$ws.Range("D5").Value = '''0.05599'
$ws.Range("G5").Value = '''13'
$ws.Range("D6").Value = '''6.486'
$ws.Range("G6").Value = '''13'
$ws.Range("D7").Value = '''3.041'
$ws.Range("G7").Value = '''13'
$ws.Range("D8").Value = '''0.8128'
$ws.Range("G8").Value = '''13'
$ws.Range("D9").Value = '''0.8408'
$ws.Range("G9").Value = '''13'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1342'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("G10").Value = '''13'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").Value = '''0.02796'
$ws.Range("E11").Value = '10BitrueCoinBTR'
$ws.Range("G11").Value = '''13'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").Value = '''0.09402'
$ws.Range("E12").Value = '11BitMartTokenBMX'
$ws.Range("G12").Value = '''13'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").Value = '''0.001526'
$ws.Range("E13").Value = '12BitForexTokenBF'
$ws.Range("G13").Value = '''13'
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = '''0.0006009'
$ws.Range("E14").Value = '13OneONEWorstin24h'
$ws.Range("G14").Value = '''13'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006108'
$ws.Range("E15").Value = '14TigerCashTCH'
$ws.Range("G15").Value = '''13'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.559'
$ws.Range("E16").Value = '15LEOLEO'
$ws.Range("G16").Value = '''13'
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").Value = '''2.118'
$ws.Range("E17").Value = '16BTSETokenBTSE'
$ws.Range("G17").Value = '''13'
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D18").Value = '''0.3188'
$ws.Range("E18").Value = '17BitpandaEcosystemTokenBEST'
$ws.Range("G18").Value = '''13'
$ws.Range("B19").Value = 'MandalaExchangeToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D19").Value = '''0.06973'
$ws.Range("E19").Value = '18MandalaExchangeTokenMDX'
$ws.Range("G19").Value = '''13'
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = '''0.03237'
$ws.Range("E20").Value = '19LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G20").Value = '''13'
$ws.Range("D21").Value = '''0.1301'
$ws.Range("G21").Value = '''13'
$ws.Range("D22").Value = '''3.736'
$ws.Range("G22").Value = '''13'
$ws.Range("D23").Value = '''0.04682'
$ws.Range("G23").Value = '''13'
$ws.Range("G24").Value = '''13'
$ws.Range("D25").Value = '''0.001248'
$ws.Range("G25").Value = '''13'
$ws.Range("D26").Value = '''0.004605'
$ws.Range("G26").Value = '''13'
$ws.Range("D27").Value = '''0.00009596'
$ws.Range("G27").Value = '''13'
$ws.Range("G28").Value = '''13'
$ws.Range("G29").Value = '''13'
$ws.Range("G30").Value = '''13'
$ws.Range("G31").Value = '''13'
$ws.Range("G32").Value = '''13'
$ws.Range("G33").Value = '''13'
$ws.Range("G34").Value = '''13'
$ws.Range("G35").Value = '''13'
$ws.Range("G36").Value = '''13'
$ws.Range("G37").Value = '''13'
$ws.Range("G38").Value = '''13'
$ws.Range("G39").Value = '''13'
$ws.Range("D40").Value = '''0.03657'
$ws.Range("G40").Value = '''13'
$ws.Range("D41").Value = '''0.006168'
$ws.Range("G41").Value = '''13'
$ws.Range("G42").Value = '''13'
$ws.Range("D43").Value = '''0.002594'
$ws.Range("G43").Value = '''13'
$ws.Range("D44").Value = '''0.008759'
$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'
$ws.Range("G44").Value = '''13'
$ws.Range("D45").Value = '''0.00005291'
$ws.Range("G45").Value = '''13'
$ws.Range("G46").Value = '''13'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("G47").Value = '''13'
$ws.Range("D48").Value = '''0.002060'
$ws.Range("G48").Value = '''13'
$ws.Range("D49").Value = '''0.00002099'
$ws.Range("G49").Value = '''13'
$ws.Range("D50").Value = '''0.0001999'
$ws.Range("G50").Value = '''13'
$ws.Range("G51").Value = '''13'
